# Updated cryptos list on Tue Jan 23 23:44:31 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text representation
# (e.g. "298.20", "39.751.75") instead of being auto-converted to numbers
# by Excel's smart cell-value parsing.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "39.751.75"
$ws.Range("E2").Value = "  +0.44%  "

$ws.Range("D3").Value = "2.238.96"
$ws.Range("E3").Value = "  -3.25%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "298.20"
$ws.Range("E5").Value = "  -2.35%  "

$ws.Range("D6").Value = "83.72"
$ws.Range("E6").Value = "  +0.39%  "

$ws.Range("D7").Value = "0.517"
$ws.Range("E7").Value = "  -1.73%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  -0.44%  "

$ws.Range("D10").Value = "30.40"
$ws.Range("E10").Value = "  +3.05%  "

$ws.Range("E11").Value = "  -2.64%  "

$ws.Range("D12").Value = "47.07"
$ws.Range("E12").Value = "  -10.27%  "

$ws.Range("E13").Value = "  -1.88%  "

$ws.Range("D14").Value = "2.588.91"
$ws.Range("E14").Value = "  -3.24%  "

$ws.Range("D15").Value = "6.34"
$ws.Range("E15").Value = "  +0.13%  "

$ws.Range("D16").Value = "14.25"
$ws.Range("E16").Value = "  -1.91%  "

$ws.Range("D17").Value = "2.241.84"

$ws.Range("D18").Value = "0.723"
$ws.Range("E18").Value = "  -2.99%  "

$ws.Range("D19").Value = "39.705.01"
$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("E20").Value = "  -1.12%  "

$ws.Range("E21").Value = "  -3.53%  "

$ws.Range("D22").Value = "65.40"
$ws.Range("E22").Value = "  -2.58%  "

$ws.Range("D23").Value = "10.51"
$ws.Range("E23").Value = "  +0.45%  "

$ws.Range("D24").Value = "229.09"
$ws.Range("E24").Value = "  -2.47%  "

$ws.Range("E25").Value = "  -0.15%  "

$ws.Range("E26").Value = "  -3.43%  "

$ws.Range("E27").Value = "  +4.59%  "

$ws.Range("D28").Value = "22.97"
$ws.Range("E28").Value = "  +0.51%  "

$ws.Range("E29").Value = "  +3.06%  "

$ws.Range("E30").Value = "  +0.56%  "

$ws.Range("D31").Value = "32.75"
$ws.Range("E31").Value = "  -3.45%  "

$ws.Range("D32").Value = "149.90"
$ws.Range("E32").Value = "  -0.09%  "

$ws.Range("E33").Value = "  -0.25%  "

# Rows 34 and 35 swap places: Filecoin <-> WEMIXToken
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "2.43"
$ws.Range("E34").Value = "  -0.73%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "4.88"
$ws.Range("E35").Value = "  -2.94%  "

$ws.Range("E36").Value = "  -0.45%  "

$ws.Range("D37").Value = "16.25"
$ws.Range("E37").Value = "  +7.02%  "

$ws.Range("E38").Value = "  -1.06%  "

$ws.Range("E39").Value = "  -0.86%  "

$ws.Range("E40").Value = "  -0.80%  "

$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("D42").Value = "3.74"
$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("D43").Value = "1.926.53"
$ws.Range("E43").Value = "  -0.09%  "

$ws.Range("D44").Value = "0.0265"
$ws.Range("E44").Value = "  +1.55%  "

$ws.Range("E45").Value = "  -8.91%  "

$ws.Range("D46").Value = "16.59"
$ws.Range("E46").Value = "  -3.67%  "

$ws.Range("D47").Value = "9.15"
$ws.Range("E47").Value = "  -1.41%  "

$ws.Range("E48").Value = "  -0.01%  "

$ws.Range("D49").Value = "2.461.15"
$ws.Range("E49").Value = "  -3.06%  "

$ws.Range("D50").Value = "71.63"
$ws.Range("E50").Value = "  +3.67%  "

$ws.Range("D51").Value = "89.21"
$ws.Range("E51").Value = "  -2.55%  "
